$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @{ "G" = 1.789499; "H" = 5.368497; "I" = 0.01244533957901722; "J" = 0.01244533957901722; "K" = 3; "L" = 1; "M" = 0.571597; "N" = 1.714791; "O" = 0.00666866079389509; "P" = 0.00666866079389509; "Q" = 1.022872259903; "R" = 9.205850339126998; "S" = 0.00008299374811730295; "T" = 0.00008299374811730295 }
  3 = @{ "G" = 1.789499; "H" = 5.368497; "I" = 0.01244533957901722; "J" = 0.01244533957901722; "K" = 3; "L" = 1; "M" = 84.83061466666666; "N" = 254.491844; "O" = 0.9896948272115175; "P" = 0.9896948272115175; "Q" = 151.8043001153853; "R" = 1366.238701038468; "S" = 0.01231708820424411; "T" = 0.01231708820424411 }
  4 = @{ "G" = 1.789499; "H" = 5.368497; "I" = 0.01244533957901722; "J" = 0.01244533957901722; "K" = 2; "L" = 0.6666666666666666; "M" = 0.3116996666666667; "N" = 0.9350989999999999; "O" = 0.00363651199458739; "P" = 0.00363651199458739; "Q" = 0.5577862418003333; "R" = 5.020076176202999; "S" = 0.00004525762665580929; "T" = 0.00004525762665580929 }
  5 = @{ "G" = 103.907654; "H" = 311.722962; "I" = 0.7226413867171911; "J" = 0.7226413867171912; "K" = 3; "L" = 1; "M" = 0.571597; "N" = 1.714791; "O" = 0.00666866079389509; "P" = 0.00666866079389509; "Q" = 59.393303303438; "R" = 534.539729730942; "S" = 0.004819050283646912; "T" = 0.004819050283646913 }
  6 = @{ "G" = 103.907654; "H" = 311.722962; "I" = 0.7226413867171911; "J" = 0.7226413867171912; "K" = 3; "L" = 1; "M" = 84.83061466666666; "N" = 254.491844; "O" = 0.9896948272115175; "P" = 0.9896948272115175; "Q" = 8814.550157391324; "R" = 79330.95141652193; "S" = 0.7151944423629618; "T" = 0.7151944423629619 }
  7 = @{ "G" = 103.907654; "H" = 311.722962; "I" = 0.7226413867171911; "J" = 0.7226413867171912; "K" = 2; "L" = 0.6666666666666666; "M" = 0.3116996666666667; "N" = 0.9350989999999999; "O" = 0.00363651199458739; "P" = 0.00363651199458739; "Q" = 32.38798111591533; "R" = 291.4918300432379; "S" = 0.00262789407058233; "T" = 0.00262789407058233 }
  8 = @{ "G" = 38.091531; "H" = 114.274593; "I" = 0.2649132737037916; "J" = 0.2649132737037916; "K" = 3; "L" = 1; "M" = 0.571597; "N" = 1.714791; "O" = 0.00666866079389509; "P" = 0.00666866079389509; "Q" = 21.773004845007; "R" = 195.957043605063; "S" = 0.001766616762130874; "T" = 0.001766616762130874 }
  9 = @{ "G" = 38.091531; "H" = 114.274593; "I" = 0.2649132737037916; "J" = 0.2649132737037916; "K" = 3; "L" = 1; "M" = 84.83061466666666; "N" = 254.491844; "O" = 0.9896948272115175; "P" = 0.9896948272115175; "Q" = 3231.327988324388; "R" = 29081.95189491949; "S" = 0.2621832966443114; "T" = 0.2621832966443114 }
  10 = @{ "G" = 38.091531; "H" = 114.274593; "I" = 0.2649132737037916; "J" = 0.2649132737037916; "K" = 2; "L" = 0.6666666666666666; "M" = 0.3116996666666667; "N" = 0.9350989999999999; "O" = 0.00363651199458739; "P" = 0.00363651199458739; "Q" = 11.873117515523; "R" = 106.858057639707; "S" = 0.0009633602973492503; "T" = 0.0009633602973492503 }
}

foreach ($r in $data.Keys) {
  foreach ($c in $data[$r].Keys) {
    $ws.Range("$c$r").Value = $data[$r][$c]
  }
}
